# Update countries & provincias Spain
#
# This script brings the "Pais" sheet's COVID-19 table up to date:
#   - 6 pairs/rotations of countries swap rank order (their row now shows
#     the other country's name because the underlying data was re-sorted
#     by total cases), so the displayed labels in column A are corrected.
#   - The numeric statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) are refreshed for
#     every row whose figures moved.
#   - The "last updated" banner in A1 is bumped to the newer timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: country labels that changed because of re-sorting ---
$ws.Range("A27").Value  = "Irak"
$ws.Range("A28").Value  = "China"
$ws.Range("A56").Value  = "Azerbaiyan"
$ws.Range("A57").Value  = "Irlanda"
$ws.Range("A58").Value  = "Ghana"
$ws.Range("A62").Value  = "Serbia"
$ws.Range("A63").Value  = "Austria"
$ws.Range("A100").Value = "Croacia"
$ws.Range("A101").Value = "Grecia"
$ws.Range("A146").Value = "Namibia"
$ws.Range("A147").Value = "Republica del Chad"
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Row 4 (Estados Unidos) ---
$ws.Range("B4").Value = 3548546
$ws.Range("C4").Value = 3469
$ws.Range("D4").Value = 1600926
$ws.Range("E4").Value = 1808431
$ws.Range("G4").Value = 46
$ws.Range("H4").Value = 139189

# --- Row 6 (India) ---
$ws.Range("D6").Value = 594723
$ws.Range("E6").Value = 322536

# --- Row 17 (Arabia Saudita) ---
$ws.Range("B17").Value = 240474
$ws.Range("C17").Value = 2671
$ws.Range("D17").Value = 183048
$ws.Range("E17").Value = 55101
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = 2325

# --- Row 24 (Canada) ---
$ws.Range("D24").Value = 47298
$ws.Range("E24").Value = 57625
$ws.Range("G24").Value = 19
$ws.Range("H24").Value = 1987

# --- Row 27 (now Irak) ---
$ws.Range("B27").Value = 83867
$ws.Range("C27").Value = 2110
$ws.Range("D27").Value = 52621
$ws.Range("E27").Value = 27814
$ws.Range("G27").Value = 87
$ws.Range("H27").Value = 3432

# --- Row 28 (now China) ---
$ws.Range("B28").Value = 83611
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 78693
$ws.Range("E28").Value = 284
$ws.Range("H28").Value = 4634

# --- Row 37 (Suecia) ---
$ws.Range("B37").Value = 56877
$ws.Range("C37").Value = 703
$ws.Range("D37").Value = 46897
$ws.Range("E37").Value = 9581
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 399

# --- Row 40 (Kuwait) ---
$ws.Range("B40").Value = 51252
$ws.Range("C40").Value = 106
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 6136

# --- Row 43 (Paises Bajos) ---
$ws.Range("B43").Value = 47426
$ws.Range("C43").Value = 375
$ws.Range("D43").Value = 32110
$ws.Range("E43").Value = 13640
$ws.Range("G43").Value = 8
$ws.Range("H43").Value = 1676

# --- Row 56 (now Azerbaiyan) ---
$ws.Range("B56").Value = 25672
$ws.Range("C56").Value = 559
$ws.Range("D56").Value = 16695
$ws.Range("E56").Value = 8651
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 326

# --- Row 57 (now Irlanda) ---
$ws.Range("B57").Value = 25670
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 23364
$ws.Range("E57").Value = 560
$ws.Range("H57").Value = 1746

# --- Row 58 (now Ghana) ---
$ws.Range("B58").Value = 25252
$ws.Range("C58").Value = 264
$ws.Range("D58").Value = 21397
$ws.Range("E58").Value = 3716
$ws.Range("H58").Value = 139

# --- Row 62 (now Serbia) ---
$ws.Range("B62").Value = 19334
$ws.Range("C62").Value = 351
$ws.Range("D62").Value = 13991
$ws.Range("E62").Value = 4914
$ws.Range("G62").Value = 11
$ws.Range("H62").Value = 429

# --- Row 63 (now Austria) ---
$ws.Range("B63").Value = 19154
$ws.Range("C63").Value = 133
$ws.Range("D63").Value = 17175
$ws.Range("E63").Value = 1269
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 710

# --- Row 67 (Uzbekistan) ---
$ws.Range("B67").Value = 14466
$ws.Range("C67").Value = 381
$ws.Range("D67").Value = 8343
$ws.Range("E67").Value = 6054
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 69

# --- Row 73 (Corea del Sur) ---
$ws.Range("B73").Value = 11252
$ws.Range("C73").Value = 461
$ws.Range("E73").Value = 8033

# --- Row 100 (now Croacia) ---
$ws.Range("B100").Value = 3953
$ws.Range("C100").Value = 92
$ws.Range("D100").Value = 2629
$ws.Range("E100").Value = 1204
$ws.Range("H100").Value = 120

# --- Row 101 (now Grecia) ---
$ws.Range("B101").Value = 3883
$ws.Range("D101").Value = 1374
$ws.Range("E101").Value = 2316
$ws.Range("H101").Value = 193

# --- Row 119 (Islandia) ---
$ws.Range("B119").Value = 1911
$ws.Range("C119").Value = 6
$ws.Range("D119").Value = 1885
$ws.Range("E119").Value = 16

# --- Row 121 (Lituania) ---
$ws.Range("C121").Value = 21

# --- Row 126 (Hong Kong) ---
$ws.Range("E126").Value = 339
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 9

# --- Row 146 (now Namibia) ---
$ws.Range("B146").Value = 960
$ws.Range("C146").Value = 96
$ws.Range("D146").Value = 31
$ws.Range("E146").Value = 927
$ws.Range("H146").Value = 2

# --- Row 147 (now Republica del Chad) ---
$ws.Range("B147").Value = 884
$ws.Range("D147").Value = 798
$ws.Range("E147").Value = 11
$ws.Range("H147").Value = 75

# --- Updated timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Julio de 2020 a las 15:26"
